# Update the PERMANOVA results table ("benthic_nmds_permanova"):
#  - Model row: Df 3 -> 2, R2 0.537 -> 0.539, F 15.45 -> 23.985
#  - Widen the 4th table column from 961 dxa (48.05pt) to 1084 dxa (54.2pt)

$d = $word.ActiveDocument
$table = $d.Tables(1)

# Row 2 ("Model") holds the statistic values that changed.
$table.Cell(2, 2).Range.Text = "2"
$table.Cell(2, 3).Range.Text = "0.539"
$table.Cell(2, 4).Range.Text = "23.985"

# Widen the 4th column (961 dxa -> 1084 dxa == 48.05pt -> 54.2pt).
$table.Columns(4).Width = 54.2
